$d = $word.ActiveDocument

# The three inline illustrations in this document are each replaced by a
# plain hyperlink run (style "Hyperlink") whose visible text is the raw
# image URL and whose target is that same URL. We always operate on
# InlineShapes.Item(1) because each deletion shifts the remaining shapes'
# indices down by one.

$replacements = @(
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Hotel/H03_Road_Buffer_and_Setbacks.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Hotel/H04_Setbacks_for_Ancillary_Structures_Substation.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Hotel/H05_Setbacks_MSCP.jpg?h=100%25&w=100%25"
)

foreach ($url in $replacements) {
    $shp = $d.InlineShapes.Item(1)
    $start = $shp.Range.Start
    $shp.Delete()
    $target = $d.Range($start, $start)
    $d.Hyperlinks.Add($target, $url, "", "", $url) | Out-Null
}

Write-Host "Replaced" $replacements.Count "illustrations with hyperlinks"
